$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1005.94446
$ws.Range("J17").Value = 1011.82855
$ws.Range("L17").Value = 3035.48565
$ws.Range("N17").Value = -3371.48565
$ws.Range("H38").Value = 597.5
$ws.Range("I38").Value = 597.5
$ws.Range("K38").Value = 1792.5
$ws.Range("M38").Value = -1420.5
$ws.Range("H87").Value = 93333.336
$ws.Range("J87").Value = 93333.336
$ws.Range("L87").Value = 93333.336
$ws.Range("N87").Value = -95829.336
$ws.Range("H90").Value = 93333.336
$ws.Range("J90").Value = 93333.336
$ws.Range("L90").Value = 280000.008
$ws.Range("N90").Value = -292480.008
$ws.Range("H112").Value = 5334.7954
$ws.Range("I112").Value = 1471
$ws.Range("K112").Value = 4413
$ws.Range("M112").Value = -3305
$ws.Range("H129").Value = 1227.8572
$ws.Range("I129").Value = 648.75
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 1946.25
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 3053.75
$ws.Range("N129").Value = -16000
$ws.Range("H132").Value = 1364.7878
$ws.Range("I132").Value = 1384
$ws.Range("K132").Value = 4152
$ws.Range("M132").Value = -1622
$ws.Range("H135").Value = 476975.1
$ws.Range("I135").Value = 625917.5600000001
$ws.Range("K135").Value = 5633258.040000001
$ws.Range("M135").Value = -5630723.040000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1547440
$ws.Range("I32").Value = 1648300.5
$ws.Range("K32").Value = 1648300.5
$ws.Range("M32").Value = -1648013.5
$ws.Range("H45").Value = 8726.357
$ws.Range("I45").Value = 2663.6667
$ws.Range("K45").Value = 2663.6667
$ws.Range("M45").Value = -2286.6667
$ws.Range("H61").Value = 7025.8125
$ws.Range("I61").Value = 3117.75
$ws.Range("J61").Value = 10933.875
$ws.Range("K61").Value = 3117.75
$ws.Range("L61").Value = 10933.875
$ws.Range("M61").Value = -2905.75
$ws.Range("N61").Value = -11357.875
$ws.Range("H74").Value = 15767.956
$ws.Range("I74").Value = 24475.424
$ws.Range("K74").Value = 24475.424
$ws.Range("M74").Value = -23601.424
$ws.Range("H77").Value = 15767.956
$ws.Range("I77").Value = 24475.424
$ws.Range("K77").Value = 122377.12
$ws.Range("M77").Value = -118009.12
$ws.Range("H101").Value = 60958.2
$ws.Range("J101").Value = 60958.2
$ws.Range("L101").Value = 60958.2
$ws.Range("N101").Value = -67448.2
$ws.Range("H110").Value = 23810768
$ws.Range("I110").Value = 1208.625
$ws.Range("K110").Value = 1208.625
$ws.Range("M110").Value = 836.375
$ws.Range("H122").Value = 11309.739
$ws.Range("I122").Value = 13669.375
$ws.Range("K122").Value = 41008.125
$ws.Range("M122").Value = -38558.125
$ws.Range("H126").Value = 5242.625
$ws.Range("I126").Value = 5242.625
$ws.Range("K126").Value = 15727.875
$ws.Range("M126").Value = -13257.875
$ws.Range("H136").Value = 7025.8125
$ws.Range("I136").Value = 3117.75
$ws.Range("J136").Value = 10933.875
$ws.Range("K136").Value = 9353.25
$ws.Range("L136").Value = 32801.625
$ws.Range("M136").Value = -6803.25
$ws.Range("N136").Value = -37901.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 4000
$ws.Range("I49").Value = 4000
$ws.Range("K49").Value = 4000
$ws.Range("M49").Value = -3761
$ws.Range("H86").Value = 28602726
$ws.Range("I86").Value = 47194.227
$ws.Range("K86").Value = 47194.227
$ws.Range("M86").Value = -46071.227
$ws.Range("H89").Value = 28602726
$ws.Range("I89").Value = 47194.227
$ws.Range("K89").Value = 235971.135
$ws.Range("M89").Value = -230355.135
$ws.Range("H94").Value = 1417.2858
$ws.Range("I94").Value = 698
$ws.Range("J94").Value = 2528.9092
$ws.Range("K94").Value = 698
$ws.Range("L94").Value = 2528.9092
$ws.Range("M94").Value = -247
$ws.Range("N94").Value = -3430.9092
$ws.Range("H99").Value = 2528394.8
$ws.Range("I99").Value = 2955
$ws.Range("K99").Value = 2955
$ws.Range("M99").Value = -1457

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3983.7273
$ws.Range("I16").Value = 2310.4167
$ws.Range("K16").Value = 2310.4167
$ws.Range("M16").Value = -2023.4167
$ws.Range("H58").Value = 13163976
$ws.Range("I58").Value = 25001368
$ws.Range("J58").Value = 11319.277
$ws.Range("K58").Value = 25001368
$ws.Range("L58").Value = 11319.277
$ws.Range("M58").Value = -25001165
$ws.Range("N58").Value = -11725.277
$ws.Range("H86").Value = 5213892
$ws.Range("I86").Value = 8933693
$ws.Range("J86").Value = 6171
$ws.Range("K86").Value = 8933693
$ws.Range("L86").Value = 6171
$ws.Range("M86").Value = -8932570
$ws.Range("N86").Value = -8417
$ws.Range("H89").Value = 5213892
$ws.Range("I89").Value = 8933693
$ws.Range("J89").Value = 6171
$ws.Range("K89").Value = 44668465
$ws.Range("L89").Value = 30855
$ws.Range("M89").Value = -44662849
$ws.Range("N89").Value = -42087
$ws.Range("H113").Value = 3983.7273
$ws.Range("I113").Value = 2310.4167
$ws.Range("K113").Value = 2310.4167
$ws.Range("M113").Value = -140.4167000000002
$ws.Range("H134").Value = 4615.415
$ws.Range("I134").Value = 1309.5862
$ws.Range("K134").Value = 3928.7586
$ws.Range("M134").Value = -1393.7586
$ws.Range("H136").Value = 13163976
$ws.Range("I136").Value = 25001368
$ws.Range("J136").Value = 11319.277
$ws.Range("K136").Value = 75004104
$ws.Range("L136").Value = 33957.831
$ws.Range("M136").Value = -75001554
$ws.Range("N136").Value = -39057.831
$ws.Range("H141").Value = 48801.383
$ws.Range("J141").Value = 51860.168
$ws.Range("L141").Value = 51860.168
$ws.Range("N141").Value = -62220.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 425.16666
$ws.Range("I23").Value = 275.25
$ws.Range("J23").Value = 725
$ws.Range("K23").Value = 825.75
$ws.Range("L23").Value = 2175
$ws.Range("M23").Value = -590.75
$ws.Range("N23").Value = -2645
$ws.Range("H33").Value = 33333700
$ws.Range("I33").Value = 66666920
$ws.Range("J33").Value = 480.6
$ws.Range("K33").Value = 400001520
$ws.Range("L33").Value = 2883.6
$ws.Range("M33").Value = -400001237
$ws.Range("N33").Value = -3449.6
$ws.Range("H39").Value = 8459.333000000001
$ws.Range("J39").Value = 9706.691999999999
$ws.Range("L39").Value = 29120.076
$ws.Range("N39").Value = -29708.076
$ws.Range("H55").Value = 27089152
$ws.Range("J55").Value = 9099352
$ws.Range("L55").Value = 27298056
$ws.Range("N55").Value = -27298410
$ws.Range("H75").Value = 2740.6667
$ws.Range("I75").Value = 1549.5
$ws.Range("J75").Value = 2978.9
$ws.Range("K75").Value = 4648.5
$ws.Range("L75").Value = 8936.700000000001
$ws.Range("M75").Value = -3650.5
$ws.Range("N75").Value = -10932.7
$ws.Range("H78").Value = 2740.6667
$ws.Range("I78").Value = 1549.5
$ws.Range("J78").Value = 2978.9
$ws.Range("K78").Value = 13945.5
$ws.Range("L78").Value = 26810.1
$ws.Range("M78").Value = -8953.5
$ws.Range("N78").Value = -36794.10000000001
$ws.Range("H103").Value = 1209.9166
$ws.Range("I103").Value = 588.3333
$ws.Range("J103").Value = 1417.1111
$ws.Range("K103").Value = 1764.9999
$ws.Range("L103").Value = 4251.3333
$ws.Range("M103").Value = -885.9999
$ws.Range("N103").Value = -6009.3333
$ws.Range("H122").Value = 1888233.5
$ws.Range("I122").Value = 4715822
$ws.Range("J122").Value = 3174.5557
$ws.Range("K122").Value = 42442398
$ws.Range("L122").Value = 28571.0013
$ws.Range("M122").Value = -42439948
$ws.Range("N122").Value = -33471.0013
$ws.Range("H134").Value = 4802.95
$ws.Range("I134").Value = 4221.3125
$ws.Range("J134").Value = 7129.5
$ws.Range("K134").Value = 12663.9375
$ws.Range("L134").Value = 21388.5
$ws.Range("M134").Value = -7593.9375
$ws.Range("N134").Value = -31528.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3611.3333
$ws.Range("I102").Value = 4136.7144
$ws.Range("K102").Value = 4136.7144
$ws.Range("M102").Value = -2514.7144
$ws.Range("H126").Value = 4292.5
$ws.Range("I126").Value = 4060.6667
$ws.Range("K126").Value = 12182.0001
$ws.Range("M126").Value = -9712.000100000001
$ws.Range("H132").Value = 6729.5293
$ws.Range("I132").Value = 1616.8334
$ws.Range("K132").Value = 4850.5002
$ws.Range("M132").Value = -2320.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4952.2915
$ws.Range("I40").Value = 4142.3076
$ws.Range("K40").Value = 4142.3076
$ws.Range("M40").Value = -4006.3076
$ws.Range("H132").Value = 13165294
$ws.Range("I132").Value = 27780398
$ws.Range("K132").Value = 83341194
$ws.Range("M132").Value = -83338664
$ws.Range("H136").Value = 7929.919
$ws.Range("I136").Value = 1788.8235
$ws.Range("K136").Value = 5366.470499999999
$ws.Range("M136").Value = -2816.470499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 3557.5
$ws.Range("I23").Value = 1373
$ws.Range("J23").Value = 10111
$ws.Range("K23").Value = 1373
$ws.Range("L23").Value = 10111
$ws.Range("M23").Value = -1144
$ws.Range("N23").Value = -10569
$ws.Range("H123").Value = 51598
$ws.Range("J123").Value = 51598
$ws.Range("L123").Value = 51598
$ws.Range("N123").Value = -61398
$ws.Range("H136").Value = 24419366
$ws.Range("I136").Value = 58824536
$ws.Range("J136").Value = 49038.082
$ws.Range("K136").Value = 176473608
$ws.Range("L136").Value = 147114.246
$ws.Range("M136").Value = -176471058
$ws.Range("N136").Value = -152214.246
